$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CDM")
$ws.Activate()

# Sessão 2 (Acompanhamento de Riscos) - preencher a linha 17 anteriormente vazia
$ws.Range("A17").Value = "Não Comprimento do Cronograma"
$ws.Range("D17").Value = "A"
$ws.Range("F17").Value = "Sim"

# Atualiza a seleção ativa da planilha para a área recém-preenchida
$ws.Range("A18:I20").Select()
